$wb = $excel.ActiveWorkbook
$gamesWs = $wb.Worksheets.Item("Games")
$nextWs = $wb.Worksheets.Item("Next")

# The game vs BOS on 45306 (previously scheduled in "Next") has now been
# played. Append it as a completed game (row 41) to the "Games" sheet with
# its final box-score stats, then remove it from the "Next" schedule sheet
# (row 2), shifting the remaining upcoming games up.

$row = 41
$gamesWs.Cells.Item($row, 1).Value = 40        # Game
$gamesWs.Cells.Item($row, 2).Value = 45306      # Date
$gamesWs.Cells.Item($row, 3).Value = -4         # Streak
$gamesWs.Cells.Item($row, 4).Value = 96         # Pts
$gamesWs.Cells.Item($row, 5).Value = 99.7       # Pace
$gamesWs.Cells.Item($row, 6).Value = 0.421      # eFG
$gamesWs.Cells.Item($row, 7).Value = 6.3        # TOV
$gamesWs.Cells.Item($row, 8).Value = 15.8       # ORB
$gamesWs.Cells.Item($row, 9).Value = 0.168      # FTR
$gamesWs.Cells.Item($row, 10).Value = 96.3      # ORT
$gamesWs.Cells.Item($row, 11).Value = "BOS"     # OppID
$gamesWs.Cells.Item($row, 12).Value = 105       # OppPts
$gamesWs.Cells.Item($row, 13).Value = 0.5       # OppeFG
$gamesWs.Cells.Item($row, 14).Value = 11.4      # OppTOV
$gamesWs.Cells.Item($row, 15).Value = 10.4      # OppORB
$gamesWs.Cells.Item($row, 16).Value = 0.28      # OppFTR
$gamesWs.Cells.Item($row, 17).Value = 105.4     # OppORT
$gamesWs.Cells.Item($row, 18).Value = 1         # Location
$gamesWs.Cells.Item($row, 19).Value = 0         # Target

# Apply the date number format (YYYY-MM-DD) used by the rest of the Date
# column, matching B2:B40.
$gamesWs.Cells.Item($row, 2).NumberFormat = $gamesWs.Cells.Item($row - 1, 2).NumberFormat

# Remove the now-played BOS game from the upcoming schedule sheet; this
# shifts every later row up by one and drops the trailing row.
$nextWs.Rows.Item(2).Delete()
